$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1 - Caribbean warming")

# Widen the columns for the new table.
$ws.Columns.Item(1).ColumnWidth = 38.71
$ws.Columns.Item(2).ColumnWidth = 19.71
$ws.Columns.Item(3).ColumnWidth = 19.71
$ws.Columns.Item(4).ColumnWidth = 19.71
$ws.Columns.Item(5).ColumnWidth = 22.71
$ws.Columns.Item(6).ColumnWidth = 22.71

# Header row.
$ws.Range("A1").Value = "Temperature.parameter"
$ws.Range("B1").Value = "HadISST..1871.2020."
$ws.Range("C1").Value = "HadISST..1981.2020."
$ws.Range("D1").Value = "HadISST..1987.2020."
$ws.Range("E1").Value = "Pathfinder..1981.2019."
$ws.Range("F1").Value = "Pathfinder..1990.2019."

# Row 2 - Caribbean Basin (degrees C per decade).
$ws.Range("A2").Value = "Caribbean Basin (°C per decade)"
$ws.Range("B2").Value = 0.04
$ws.Range("C2").Value = 0.17
$ws.Range("D2").Value = 0.18
$ws.Range("E2").Value = 0.17
$ws.Range("F2").Value = "NA"

# Row 3 - Caribbean Basin (total degrees C for period).
$ws.Range("A3").Value = "Caribbean Basin (total °C for period)"
$ws.Range("B3").Value = 0.6
$ws.Range("C3").Value = 0.68
$ws.Range("D3").Value = 0.61
$ws.Range("E3").Value = 0.66
$ws.Range("F3").Value = "NA"

# Row 4 - Caribbean Reefs (degrees C per decade).
$ws.Range("A4").Value = "Caribbean Reefs (°C per decade)"
$ws.Range("B4").Value = 0.04
$ws.Range("C4").Value = 0.15
$ws.Range("D4").Value = 0.16
$ws.Range("E4").Value = 0.19
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "0.2"
$ws.Range("F4").Style = "Normal"

# Row 5 - Caribbean Reefs (total degrees C for period).
$ws.Range("A5").Value = "Caribbean Reefs (total °C for period)"
$ws.Range("B5").Value = 0.6
$ws.Range("C5").Value = 0.6
$ws.Range("D5").Value = 0.54
$ws.Range("E5").Value = 0.74
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0.66"
$ws.Range("F5").Style = "Normal"
